# RPA datasets push 2024-06-19
# Insert a new IPO record (디비금융스팩12호 / DB) as row 3, shifting the
# existing rows (old row 3 onward) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 3 (old row 2's data stays,
# everything from old row 3 down shifts to row 4+).
$ws.Rows.Item(3).Insert()

# The date-like columns (A, O, P) must stay plain text (as the rest of the
# sheet stores them), so force Text format before writing - otherwise Excel
# auto-converts "2024-06-18"-style strings into real date serials.
$ws.Range("A3:A3").NumberFormat = "@"
$ws.Range("O3:O3").NumberFormat = "@"
$ws.Range("P3:P3").NumberFormat = "@"

# Populate the newly inserted row 3 with the new listing's data.
$ws.Range("A3").Value = "2024-06-18"
$ws.Range("B3").Value = "디비금융스팩12호"
$ws.Range("C3").Value = "코스닥"
$ws.Range("D3").Value = 100
$ws.Range("E3").Value = "DB"
$ws.Range("F3").Value = 100
$ws.Range("G3").Value = "-"
$ws.Range("H3").Value = "-"
$ws.Range("I3").Value = "-"
$ws.Range("J3").Value = "-"
$ws.Range("K3").Value = "대표"
$ws.Range("L3").Value = "-"
$ws.Range("M3").Value = 2000
$ws.Range("N3").Value = 100
$ws.Range("O3").Value = "2024-06-05"
$ws.Range("P3").Value = "2024-06-11"
$ws.Range("Q3").Value = 3750000
